$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header-like row 2 with text values (will create sharedStrings.xml)
$ws.Range("A2").Value = "b"
$ws.Range("B2").Value = "AB"

# Flip sign of B3, B4, B5 (now positive)
$ws.Range("B3").Value = 89.3
$ws.Range("B4").Value = 84.7
$ws.Range("B5").Value = 76.599999999999994

# Update formula for C3 and the shared formula C4:C10 so it subtracts A from B
$ws.Range("C3").Formula = "=(B3-A3)/31.75"
$ws.Range("C4:C10").Formula = "=(B4-A4)/31.75"

# Update the selected cell/range to D10, matching the new view state
$ws.Range("D10").Select() | Out-Null
